$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NhapSanPham")

# Update the data cells that were re-imported with corrected/renumbered values
$ws.Range("B2").Value = "sai rồi nè 1"
$ws.Range("C3").Value = "sai nữa nè 2"
$ws.Range("A4").Value = "Phân bón fail fail 3"
$ws.Range("B4").Value = "sai rồi nè 3"
$ws.Range("A5").Value = "Phân bón fail fail 4"
$ws.Range("C5").Value = "sai nữa nè 4"

# Update the active selection to match the author's final cursor position
$ws.Activate()
$ws.Range("A5").Select()
